# Add a new "2023" column (U) to the right of the existing "2022" column (T),
# copying the formatting from column T and filling in the 2023 data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the "2022" column (T4:T14) into the new
# "2023" column (U4:U14); this also extends dimension/row spans automatically.
$ws.Range("T4:T14").Copy($ws.Range("U4:U14"))

# Header
$ws.Range("U4").Value = 2023

# Data values for 2023
$ws.Range("U5").Value = 0.5
$ws.Range("U6").Value = 0.3
$ws.Range("U7").Value = 0.4
$ws.Range("U8").Value = 0.4
$ws.Range("U9").Value = 3.2
$ws.Range("U10").Value = 0.6
$ws.Range("U11").Value = "-"
$ws.Range("U12").Value = 0.6
$ws.Range("U13").Value = 0.1
$ws.Range("U14").Value = 0.5

# Update the sheet's selection to match the target workbook (B1 instead of V7)
$ws.Range("B1").Select() | Out-Null
